# spring 24 week 12 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.31
$ws.Range("D2").Value = 1.27

$ws.Range("B3").Value = 1.53

$ws.Range("B4").Value = 1.45
$ws.Range("E4").Value = 1.22
$ws.Range("G4").Value = 1

$ws.Range("D5").Value = 1.35
$ws.Range("E5").Value = 1.22

$ws.Range("G6").Value = 1

$ws.Range("D7").Value = 1.73
